# Append new water-usage rows to the "Data" sheet (rows 21-37).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$dateLabel = "21/8/2024"

$values = @(
    "519 liters",
    "519 liters",
    "519 liters",
    "519 liters",
    "520 liters",
    "521 liters",
    "522 liters",
    "523 liters",
    "524 liters",
    "525 liters",
    "526 liters",
    "527 liters",
    "528 liters",
    "529 liters",
    "530 liters",
    "530 liters",
    "531 liters"
)

$startRow = 21
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dateLabel
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
